$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1) - text is unchanged ("Task Number", "Interval (Days)",
# "Men", "Menhours") so nothing to edit there.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Data rows 2-7: replace task numbers / interval values (old rows reused)
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,1).Value = "524504-50-01"
$ws.Cells.Item(2,2).Value = 30
$ws.Cells.Item(2,3).Value = 1
$ws.Cells.Item(2,4).Value = 0.2

$ws.Cells.Item(3,1).Value = "212700-00-01"
$ws.Cells.Item(3,2).Value = 365
$ws.Cells.Item(3,3).Value = 1
$ws.Cells.Item(3,4).Value = 0.3

$ws.Cells.Item(4,1).Value = "262400-00-02"
$ws.Cells.Item(4,2).Value = 365
$ws.Cells.Item(4,3).Value = 1
$ws.Cells.Item(4,4).Value = 0.3

$ws.Cells.Item(5,1).Value = "335000-02-04"
$ws.Cells.Item(5,2).Value = 365
$ws.Cells.Item(5,3).Value = 1
$ws.Cells.Item(5,4).Value = 1

$ws.Cells.Item(6,1).Value = "342200-00-01"
$ws.Cells.Item(6,2).Value = 365
$ws.Cells.Item(6,3).Value = 1
$ws.Cells.Item(6,4).Value = 13

$ws.Cells.Item(7,1).Value = "341100-50-01"
$ws.Cells.Item(7,2).Value = 365
$ws.Cells.Item(7,3).Value = 1
$ws.Cells.Item(7,4).Value = 2

# ---------------------------------------------------------------------------
# Remove the old per-cell styling (borders / wrap / right-align) on rows 2-7
# so they fall back to the default (unstyled) cell format, matching the new
# plain data rows.
# ---------------------------------------------------------------------------
$ws.Range("A2:D7").Style = "Normal"

# AutoFit clears the custom row height + thick-bottom-border hint that used
# to come from the removed styles; re-apply the literal heights that survive
# in the target for rows 2/3/5/6 (rows 4 and 7 end up with no explicit
# height at all).
for ($r = 2; $r -le 7; $r++) {
    $ws.Rows.Item($r).AutoFit()
}
$ws.Rows.Item(2).RowHeight = 14.4
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 14.4
$ws.Rows.Item(6).RowHeight = 15

# ---------------------------------------------------------------------------
# New data rows 8-19
# ---------------------------------------------------------------------------
$ws.Cells.Item(8,1).Value = "494000-00-01"
$ws.Range("B8").Formula = "=ROUNDDOWN(1000/14,0)"
$ws.Cells.Item(8,3).Value = 1
$ws.Cells.Item(8,4).Value = 0.1

$ws.Cells.Item(9,1).Value = "521102-00-09"
$ws.Cells.Item(9,2).Value = 1000
$ws.Cells.Item(9,3).Value = 2
$ws.Cells.Item(9,4).Value = 0.2

$ws.Cells.Item(10,1).Value = "801000-00-02"
$ws.Cells.Item(10,2).Value = 1000
$ws.Cells.Item(10,3).Value = 1
$ws.Cells.Item(10,4).Value = 0.2

$ws.Cells.Item(11,1).Value = "801000-00-02"
$ws.Cells.Item(11,2).Value = 1000
$ws.Cells.Item(11,3).Value = 1
$ws.Cells.Item(11,4).Value = 0.3

$ws.Cells.Item(12,1).Value = "531019-00-05"
$ws.Cells.Item(12,2).Value = 1000
$ws.Cells.Item(12,3).Value = 1
$ws.Cells.Item(12,4).Value = 0.3

$ws.Cells.Item(13,1).Value = "062403-00-01"
$ws.Range("B13").Formula = "=ROUNDDOWN(1000/13,0)"
$ws.Cells.Item(13,3).Value = 1
$ws.Cells.Item(13,4).Value = 0.3

$ws.Cells.Item(14,1).Value = "062403-00-01"
$ws.Cells.Item(14,3).Value = 1
$ws.Cells.Item(14,4).Value = 0.5

$ws.Cells.Item(15,1).Value = "215000-00-09"
$ws.Cells.Item(15,3).Value = 1
$ws.Cells.Item(15,4).Value = 0.7

$ws.Cells.Item(16,1).Value = "241000-00-02"
$ws.Cells.Item(16,3).Value = 1
$ws.Cells.Item(16,4).Value = 0.8

$ws.Cells.Item(17,1).Value = "241000-00-02"
$ws.Cells.Item(17,3).Value = 1
$ws.Cells.Item(17,4).Value = 0.1

$ws.Cells.Item(18,1).Value = "254201-00-01"
$ws.Cells.Item(18,3).Value = 1
$ws.Cells.Item(18,4).Value = 0.1

# Shared formula block for B14:B18
$ws.Range("B14:B18").Formula = "=ROUNDDOWN(1000/13,0)"

$ws.Cells.Item(19,1).Value = "255100-00-01"
$ws.Range("B19").Formula = "=ROUNDDOWN(1000/13,0)"
$ws.Cells.Item(19,3).Value = 1
$ws.Cells.Item(19,4).Value = 0.3

# ---------------------------------------------------------------------------
# Selection
# ---------------------------------------------------------------------------
$null = $ws.Range("I16").Select()
